# Remove the trailing "Ver no Jupiter ..." blurb and the copyright/footer
# paragraph (plus the blank paragraph that separated them from the
# "LOB1018: Física I (Requisito)" requirement line), while leaving the
# "LOB1018..." paragraph and the final blank/page-break paragraphs intact.

$d = $word.ActiveDocument

$start = $null
$end = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $text = $d.Paragraphs.Item($i).Range.Text

    if ($text -like "*LOB1018*F*sica I (Requisito)*") {
        # The blank paragraph right after this one starts the block to delete.
        $start = $d.Paragraphs.Item($i + 1).Range.Start
    }
    if ($text -like "*Powered by Jekyll*") {
        $end = $d.Paragraphs.Item($i).Range.End
    }
}

if ($start -ne $null -and $end -ne $null) {
    $r = $d.Range($start, $end)
    $r.Delete()
}
